$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.735.48"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "2.290.36"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'268.78"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").Value = "'93.49"
$ws.Range("E6").Value = "  +8.69%  "
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "  +1.55%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.618"
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("D10").Value = "'45.28"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("D12").Value = "'8.03"
$ws.Range("E12").Value = "  +5.43%  "
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "2.633.50"
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("D15").Value = "'15.23"
$ws.Range("E15").Value = "  +4.03%  "
$ws.Range("D16").Value = "'0.848"
$ws.Range("E16").Value = "  +8.38%  "
$ws.Range("D17").Value = "2.309.12"
$ws.Range("E17").Value = "  +4.44%  "
$ws.Range("D18").Value = "43.682.42"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").Value = "'6.24"
$ws.Range("E20").Value = "  +4.07%  "
$ws.Range("D21").Value = "'70.99"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").Value = "'2.28"
$ws.Range("E22").Value = "  -4.94%  "
$ws.Range("D23").Value = "'235.98"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").Value = "'9.66"
$ws.Range("E24").Value = "  +6.60%  "
$ws.Range("D25").Value = "'0.998"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "'2.50"
$ws.Range("E26").Value = "  +10.03%  "
$ws.Range("D27").Value = "'11.20"
$ws.Range("E27").Value = "  +3.45%  "
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("D29").Value = "'39.14"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").Value = "'22.15"
$ws.Range("E31").Value = "  +7.80%  "
$ws.Range("D32").Value = "'173.05"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").Value = "'0.0881"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").Value = "'5.53"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("D36").Value = "'0.110"
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("D37").Value = "'4.51"
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("E39").Value = "  +3.14%  "
$ws.Range("D40").Value = "'0.234"
$ws.Range("E40").Value = "  +14.35%  "
$ws.Range("D41").Value = "'2.30"
$ws.Range("E41").Value = "  +9.51%  "
$ws.Range("D42").Value = "'12.23"
$ws.Range("E42").Value = "  -1.58%  "
$ws.Range("D43").Value = "'1.31"
$ws.Range("E43").Value = "  +15.77%  "
$ws.Range("D44").Value = "'5.44"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "'61.09"
$ws.Range("E45").Value = "  -5.73%  "
$ws.Range("E46").Value = "  +5.92%  "
$ws.Range("D47").Value = "'0.101"
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("D48").Value = "'99.79"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").Value = "2.511.46"
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("E51").Value = "  -4.50%  "
